$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -7
$ws.Range("F6").Value = -8
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = -12
$ws.Range("F9").Value = -2
$ws.Range("F11").Value = 5
$ws.Range("F17").Value = -5
$ws.Range("F18").Value = -4
